$d = $word.ActiveDocument

# Create the three new paragraph styles, each based on MSCParagraph,
# matching the added <w:style> entries in styles.xml.
$newStyleNames = @("MSC_Paragraph_A", "MSC_Paragraph_B", "MSC_Paragraph_C")
foreach ($name in $newStyleNames) {
    $newStyle = $d.Styles.Add($name, 1)
    $newStyle.BaseStyle = $d.Styles("MSCParagraph")
}

$paras = $d.Paragraphs

# Paragraphs that previously used "MSCParagraph" now use "MSC_Paragraph_A".
$paragraphAIndexes = @(8, 12, 23, 27)
foreach ($idx in $paragraphAIndexes) {
    $paras.Item($idx).Style = "MSC_Paragraph_A"
}

# Paragraphs in column B of the two scripture text tables now use "MSC_Paragraph_B".
$paragraphBIndexes = @(13, 14, 15, 28, 29, 30)
foreach ($idx in $paragraphBIndexes) {
    $paras.Item($idx).Style = "MSC_Paragraph_B"
}

# Paragraphs in column C of the two scripture text tables now use "MSC_Paragraph_C".
$paragraphCIndexes = @(16, 17, 18, 31, 32, 33)
foreach ($idx in $paragraphCIndexes) {
    $paras.Item($idx).Style = "MSC_Paragraph_C"
}

# Copyright paragraphs in columns B and C now explicitly use "MSCCopyright".
$copyrightIndexes = @(37, 38)
foreach ($idx in $copyrightIndexes) {
    $paras.Item($idx).Style = "MSCCopyright"
}
